# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update the "K" column (G) values for rows 2-35 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 1
    6  = 0
    7  = 2
    8  = 1
    9  = 3
    10 = 2
    11 = 0
    12 = 2
    13 = 1
    14 = 0
    15 = 2
    16 = 3
    17 = 2
    18 = 0
    19 = 1
    20 = 0
    21 = 3
    22 = 1
    23 = 1
    24 = 1
    25 = 4
    26 = 1
    27 = 0
    28 = 0
    29 = 1
    30 = 2
    31 = 3
    32 = 1
    33 = 2
    34 = 1
    35 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
